# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables header rows encode pseudo-attributes inside plain text
# cells (e.g. "!!!ObjTables ObjTablesVersion='0.0.8'"). This change
# renames those pseudo-attribute keys to lowerCamelCase:
#   ObjTablesVersion -> objTablesVersion
#   Type             -> type
#   Id               -> id

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws1.Range("A2").Value = "!!ObjTables type='Data' id='SimpleModel'"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "!!ObjTables type='Data' id='ExtraSheet'"
